$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 3.8
$ws.Range("G2").Value = 4.4
$ws.Range("H2").Value = 2.06
$ws.Range("J2").Value = 3.3
$ws.Range("K2").Value = 3.75
$ws.Range("P2").Value = 1.76
$ws.Range("T2").Value = 1.87
$ws.Range("W2").Value = 1.3
$ws.Range("AG3").Value = 990
$ws.Range("F3").Value = 1.88
$ws.Range("G3").Value = 1.93
$ws.Range("J3").Value = 3.75
$ws.Range("K3").Value = 3.9
$ws.Range("U3").Value = 1.96
$ws.Range("W3").Value = 2.08
$ws.Range("AB4").Value = 12
$ws.Range("AG4").Value = 13.5
$ws.Range("F4").Value = 2.64
$ws.Range("G4").Value = 2.98
$ws.Range("H4").Value = 2.8
$ws.Range("I4").Value = 3.15
$ws.Range("J4").Value = 3.1
$ws.Range("K4").Value = 3.55
$ws.Range("L4").Value = 1.49
$ws.Range("N4").Value = 3.15
$ws.Range("Q4").Value = 2.2
$ws.Range("S4").Value = 4
$ws.Range("U4").Value = 2
$ws.Range("V4").Value = 1.47
$ws.Range("W4").Value = 1.51
$ws.Range("AB5").Value = 8.199999999999999
$ws.Range("AC5").Value = 10.5
$ws.Range("AD5").Value = 36
$ws.Range("AG5").Value = 10.5
$ws.Range("AJ5").Value = 17
$ws.Range("AK5").Value = 20
$ws.Range("AN5").Value = 8.800000000000001
$ws.Range("G5").Value = 1.57
$ws.Range("H5").Value = 7
$ws.Range("I5").Value = 8.4
$ws.Range("L5").Value = 1.39
$ws.Range("N5").Value = 3.9
$ws.Range("O5").Value = 1.28
$ws.Range("P5").Value = 2.02
$ws.Range("Q5").Value = 1.84
$ws.Range("R5").Value = 1.39
$ws.Range("S5").Value = 3.1
$ws.Range("T5").Value = 1.94
$ws.Range("U5").Value = 1.88
$ws.Range("W5").Value = 2.66
$ws.Range("X5").Value = 19
$ws.Range("G6").Value = 3.35
$ws.Range("H6").Value = 2.34
$ws.Range("I6").Value = 2.44
$ws.Range("J6").Value = 3.6
$ws.Range("N6").Value = 3.7
$ws.Range("O6").Value = 1.32
$ws.Range("Q6").Value = 2.02
$ws.Range("R6").Value = 1.34
$ws.Range("S6").Value = 3.6
$ws.Range("T6").Value = 1.75
$ws.Range("U6").Value = 2.14
$ws.Range("V6").Value = 1.69
$ws.Range("W6").Value = 1.43
$ws.Range("AB7").Value = 11
$ws.Range("AC7").Value = 14.5
$ws.Range("AF7").Value = 9
$ws.Range("AG7").Value = 12
$ws.Range("AJ7").Value = 12
$ws.Range("AK7").Value = 14.5
$ws.Range("AN7").Value = 5.6
$ws.Range("F7").Value = 1.37
$ws.Range("G7").Value = 1.41
$ws.Range("H7").Value = 9.4
$ws.Range("I7").Value = 11
$ws.Range("J7").Value = 5.3
$ws.Range("K7").Value = 5.8
$ws.Range("L7").Value = 1.32
$ws.Range("N7").Value = 5
$ws.Range("P7").Value = 2.36
$ws.Range("Q7").Value = 1.67
$ws.Range("R7").Value = 1.54
$ws.Range("S7").Value = 2.7
$ws.Range("T7").Value = 1.92
$ws.Range("U7").Value = 1.94
$ws.Range("W7").Value = 3.4
$ws.Range("Y7").Value = 42
$ws.Range("AB8").Value = 6.8
$ws.Range("AC8").Value = 8.4
$ws.Range("F8").Value = 1.8
$ws.Range("G8").Value = 1.91
$ws.Range("H8").Value = 5.7
$ws.Range("I8").Value = 6.8
$ws.Range("J8").Value = 3.25
$ws.Range("K8").Value = 3.6
$ws.Range("L8").Value = 1.59
$ws.Range("M8").Value = 1.13
$ws.Range("N8").Value = 2.58
$ws.Range("P8").Value = 1.5
$ws.Range("Q8").Value = 2.72
$ws.Range("S8").Value = 5.6
$ws.Range("T8").Value = 2.28
$ws.Range("U8").Value = 1.62
$ws.Range("V8").Value = 1.18
$ws.Range("W8").Value = 2.1
$ws.Range("X8").Value = 55
$ws.Range("F9").Value = 1.99
$ws.Range("G9").Value = 2.22
$ws.Range("H9").Value = 3.45
$ws.Range("I9").Value = 4.4
$ws.Range("J9").Value = 3.3
$ws.Range("O9").Value = 1.28
$ws.Range("T9").Value = 1.67
$ws.Range("U9").Value = 2.1
$ws.Range("W9").Value = 1.83
